$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1:H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(6,7),
    @(5,6),
    @(3,5),
    @(5,6),
    @(6,7),
    @(7,8),
    @(6,8),
    @(5,7),
    @(6,9),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,4),
    @(1,6),
    @(1,3),
    @(4,6),
    @(4,5),
    @(5,7),
    @(7,8),
    @(8,9),
    @(2,4),
    @(4,9),
    @(7,7),
    @(1,3),
    @(1,4),
    @(5,7),
    @(8,9),
    @(3,4),
    @(1,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
